# #CRM-258 add extra columns in buyback advance search download
#
# Adds four new columns (M:P) to the BuybackOrderSnapshot template:
#   SF Tax / {order:cp_tax_charge}
#   GST Amount / {order:gst_amount}
#   Partner Sweetner Charges / {order:partner_sweetner_charges}
#   Claimed Price / {order:cp_claimed_price}
#
# Row 1 holds the bold column headers, row 2 holds the merge-field
# placeholders, matching the formatting already used for the existing
# columns (A:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - bold, like the rest of the header row
$ws.Range("M1").Value = "SF Tax"
$ws.Range("N1").Value = "GST Amount"
$ws.Range("O1").Value = "Partner Sweetner Charges"
$ws.Range("P1").Value = "Claimed Price"
$ws.Range("M1:P1").Font.Bold = $true

# New placeholder cells (row 2) - regular weight, left aligned like the
# rest of the data row
$ws.Range("M2").Value = "{order:cp_tax_charge}"
$ws.Range("N2").Value = "{order:gst_amount}"
$ws.Range("O2").Value = "{order:partner_sweetner_charges}"
$ws.Range("P2").Value = "{order:cp_claimed_price}"
$ws.Range("M2:P2").Font.Bold = $false
$ws.Range("M2:P2").HorizontalAlignment = -4131
